$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2").Value = "2025-11-13T06:52:14.785132"
$ws.Range("Z3").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z4").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z5").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z6").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z7").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z8").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z9").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z10").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z11").Value = "2025-11-13T06:52:14.786130"
$ws.Range("Z12").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z13").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z14").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z15").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z16").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z17").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z18").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z19").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z20").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z21").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z22").Value = "2025-11-13T06:52:14.787130"
$ws.Range("Z23").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z24").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z25").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z26").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z27").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z28").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z29").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z30").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z31").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z32").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z33").Value = "2025-11-13T06:52:14.788130"
$ws.Range("Z34").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z35").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z36").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z37").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z38").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z39").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z40").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z41").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z42").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z43").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z44").Value = "2025-11-13T06:52:14.789129"
$ws.Range("Z45").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z46").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z47").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z48").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z49").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z50").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z51").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z52").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z53").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z54").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z55").Value = "2025-11-13T06:52:14.790129"
$ws.Range("Z56").Value = "2025-11-13T06:52:14.791129"
$ws.Range("Z57").Value = "2025-11-13T06:52:14.791129"
$ws.Range("Z58").Value = "2025-11-13T06:52:14.897646"
$ws.Range("Z59").Value = "2025-11-13T06:52:14.897646"
$ws.Range("Z60").Value = "2025-11-13T06:52:14.897646"
$ws.Range("Z61").Value = "2025-11-13T06:52:14.898646"
$ws.Range("Z62").Value = "2025-11-13T06:52:14.898646"
$ws.Range("Z63").Value = "2025-11-13T06:52:14.898646"
$ws.Range("Z64").Value = "2025-11-13T06:52:14.898646"
$ws.Range("Z65").Value = "2025-11-13T06:52:14.899645"
$ws.Range("Z66").Value = "2025-11-13T06:52:14.899645"
$ws.Range("Z67").Value = "2025-11-13T06:52:14.899645"
$ws.Range("Z68").Value = "2025-11-13T06:52:14.899645"
$ws.Range("Z69").Value = "2025-11-13T06:52:14.900645"
$ws.Range("Z70").Value = "2025-11-13T06:52:14.900645"
$ws.Range("Z71").Value = "2025-11-13T06:52:15.066129"
$ws.Range("Z72").Value = "2025-11-13T06:52:15.066129"
$ws.Range("Z73").Value = "2025-11-13T06:52:15.066129"
$ws.Range("Z74").Value = "2025-11-13T06:52:15.066129"
$ws.Range("Z75").Value = "2025-11-13T06:52:15.066129"
$ws.Range("Z76").Value = "2025-11-13T06:52:15.066129"
$ws.Range("Z77").Value = "2025-11-13T06:52:15.066129"
$ws.Range("Z78").Value = "2025-11-13T06:52:15.067130"
$ws.Range("Z79").Value = "2025-11-13T06:52:15.067130"
